$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overal Stats")
$ws.Range("C1").Value = "test"
Write-Host "Value set"
